# Fix typo: "Benoulli" -> "Bernoulli" (occurs 4 times in the document)
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("Benoulli", $false, $false, $false, $false, $false, $true, 1, $false, "Bernoulli", 2)
